# Update title and licensing slides for ATPESC
#
# The presentation's second slide ("License, Citation and Acknowledgements")
# has its tutorial-citation paragraph rewritten: the citation drops
# "Patricia A. Grubel" from the author list, swaps the venue/track name
# from "Better Scientific Software tutorial, in ISC High Performance" to
# "Software Productivity and Sustainability track, in Argonne Training
# Program on Extreme-Scale Computing (ATPESC)", and the cited DOI number
# is updated to the new figshare record.

$p = $ppt.ActivePresentation

# Slide 2 = "License, Citation and Acknowledgements"
$s = $p.Slides.Item(2)

# Shape 2 = "Content Placeholder 2", holding the License/Citation/
# Acknowledgements body text.
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Update the citation sentence (run immediately before the DOI link) ---
$oldCitation = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Patricia A. Grubel, Rinku K. Gupta, and David M. Rogers, Better Scientific Software tutorial, in ISC High Performance, online, 2021. DOI: "
$newCitation = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Rinku K. Gupta, and David M. Rogers, Software Productivity and Sustainability track, in Argonne Training Program on Extreme-Scale Computing (ATPESC), online, 2021. DOI: "

$found = $tr.Find($oldCitation)
if ($found -ne $null) {
    $found.Text = $newCitation
}

# --- Update the cited DOI number (the hyperlinked run right after it) ---
$oldDoi = "10.6084/m9.figshare.14642520"
$newDoi = "10.6084/m9.figshare.15130590"

$found2 = $tr.Find($oldDoi)
if ($found2 -ne $null) {
    $found2.Text = $newDoi
}
